$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# row 10
$ws.Range("H10").Value = 5313
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 10526
$ws.Range("K10").Value = 100
$ws.Range("L10").Value = 10526
$ws.Range("M10").Value = 193
$ws.Range("N10").Value = -11112
# row 11
$ws.Range("H11").Value = 21.428572
$ws.Range("I11").Value = 21.428572
$ws.Range("K11").Value = 21.428572
$ws.Range("M11").Value = 118.571428
# row 17
$ws.Range("H17").Value = 2452.1428
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2694.1667
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 8082.500100000001
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -8418.500100000001
# row 21
$ws.Range("H21").Value = 2593.3333
$ws.Range("I21").Value = 3481
$ws.Range("J21").Value = 1883.2
$ws.Range("K21").Value = 3481
$ws.Range("L21").Value = 1883.2
$ws.Range("M21").Value = -3013
$ws.Range("N21").Value = -2819.2
# row 23
$ws.Range("H23").Value = 2593.3333
$ws.Range("I23").Value = 3481
$ws.Range("J23").Value = 1883.2
$ws.Range("K23").Value = 3481
$ws.Range("L23").Value = 1883.2
$ws.Range("M23").Value = -3247
$ws.Range("N23").Value = -2351.2
# row 40
$ws.Range("H40").Value = 4824.227
$ws.Range("J40").Value = 6967.5713
$ws.Range("L40").Value = 6967.5713
$ws.Range("N40").Value = -7317.5713
# row 41
$ws.Range("H41").Value = 2518.5
$ws.Range("I41").Value = 2499.5
$ws.Range("J41").Value = 2524.8333
$ws.Range("K41").Value = 2499.5
$ws.Range("L41").Value = 2524.8333
$ws.Range("M41").Value = -2059.5
$ws.Range("N41").Value = -3404.8333
# row 43
$ws.Range("H43").Value = 2283.3333
$ws.Range("J43").Value = 1600
$ws.Range("L43").Value = 1600
$ws.Range("N43").Value = -1738
# row 53
$ws.Range("H53").Value = 431.33334
$ws.Range("I53").Value = 668
$ws.Range("K53").Value = 668
$ws.Range("M53").Value = -31
# row 62
$ws.Range("H62").Value = 7859.4165
$ws.Range("I62").Value = 3759
$ws.Range("K62").Value = 3759
$ws.Range("M62").Value = -3135
# row 64
$ws.Range("H64").Value = 6543.8887
$ws.Range("J64").Value = 8333.333000000001
$ws.Range("L64").Value = 8333.333000000001
$ws.Range("N64").Value = -8829.333000000001
# row 65
$ws.Range("H65").Value = 7859.4165
$ws.Range("I65").Value = 3759
$ws.Range("K65").Value = 18795
$ws.Range("M65").Value = -15675
# row 67
$ws.Range("H67").Value = 6543.8887
$ws.Range("J67").Value = 8333.333000000001
$ws.Range("L67").Value = 8333.333000000001
$ws.Range("N67").Value = -10049.333
# row 80
$ws.Range("H80").Value = 905.1667
$ws.Range("I80").Value = 357.75
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1073.25
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -75.25
$ws.Range("N80").Value = -7996
# row 83
$ws.Range("H83").Value = 905.1667
$ws.Range("I83").Value = 357.75
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 3219.75
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = 1772.25
$ws.Range("N83").Value = -27984
# row 86
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = ""
$ws.Range("N86").Value = 0
$ws.Range("M86").Value = -877
# row 89
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = ""
$ws.Range("N89").Value = 0
$ws.Range("M89").Value = -4384
# row 104
$ws.Range("H104").Value = 1040.25
$ws.Range("I104").Value = 1040.25
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 3120.75
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = ""
$ws.Range("N104").Value = -1373.75
# row 106
$ws.Range("H106").Value = 4466.6665
$ws.Range("J106").Value = 4200
$ws.Range("L106").Value = 4200
$ws.Range("N106").Value = -5462
# row 112
$ws.Range("H112").Value = 3066
$ws.Range("J112").Value = 3599
$ws.Range("L112").Value = 10797
$ws.Range("N112").Value = -13013
# row 113
$ws.Range("H113").Value = 2285.7144
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 1054
$ws.Range("N113").Value = -9008
# row 114
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").Value = ""
# row 116
$ws.Range("H116").Value = 7954.8
$ws.Range("J116").Value = 4499.6665
$ws.Range("L116").Value = 4499.6665
$ws.Range("N116").Value = -11383.6665
# row 129
$ws.Range("H129").Value = 2242.25
$ws.Range("I129").Value = 1489.5
$ws.Range("J129").Value = 2493.1667
$ws.Range("K129").Value = 4468.5
$ws.Range("L129").Value = 7479.500100000001
$ws.Range("M129").Value = 531.5
$ws.Range("N129").Value = -17479.5001
# row 134
$ws.Range("H134").Value = 191853.33
$ws.Range("J134").Value = 191853.33
$ws.Range("L134").Value = 191853.33
$ws.Range("N134").Value = -201993.33
# row 137
$ws.Range("H137").Value = 2955.5454
$ws.Range("I137").Value = 1377.75
$ws.Range("J137").Value = 3857.1428
$ws.Range("K137").Value = 4133.25
$ws.Range("L137").Value = 11571.4284
$ws.Range("M137").Value = -1583.25
$ws.Range("N137").Value = -16671.4284
# row 138
$ws.Range("H138").Value = 2612.6667
$ws.Range("J138").Value = 5333.3335
$ws.Range("L138").Value = 16000.0005
$ws.Range("N138").Value = -26280.0005

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 3550.5
$ws.Range("I2").Value = 2889.5
$ws.Range("K2").Value = 2889.5
$ws.Range("M2").Value = -2776.5
# row 32
$ws.Range("H32").Value = 4406.839
$ws.Range("I32").Value = 3620.4
$ws.Range("J32").Value = 28000
$ws.Range("K32").Value = 3620.4
$ws.Range("L32").Value = 28000
$ws.Range("M32").Value = -3333.4
$ws.Range("N32").Value = -28574
# row 45
$ws.Range("H45").Value = 2880.4
$ws.Range("I45").Value = 2100.5
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 2100.5
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -1723.5
$ws.Range("N45").Value = -6754
# row 61
$ws.Range("H61").Value = 2864
$ws.Range("I61").Value = 2436.8
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2436.8
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2224.8
$ws.Range("N61").Value = -5424
# row 74
$ws.Range("H74").Value = 8698.85
$ws.Range("I74").Value = 8666.4375
$ws.Range("J74").Value = 8828.5
$ws.Range("K74").Value = 8666.4375
$ws.Range("L74").Value = 8828.5
$ws.Range("M74").Value = -7792.4375
$ws.Range("N74").Value = -10576.5
# row 77
$ws.Range("H77").Value = 8698.85
$ws.Range("I77").Value = 8666.4375
$ws.Range("J77").Value = 8828.5
$ws.Range("K77").Value = 43332.1875
$ws.Range("L77").Value = 44142.5
$ws.Range("M77").Value = -38964.1875
$ws.Range("N77").Value = -52878.5
# row 96
$ws.Range("H96").Value = 2534352
$ws.Range("J96").Value = 2534352
$ws.Range("L96").Value = 2534352
$ws.Range("N96").Value = -2539844
# row 98
$ws.Range("H98").Value = 15000
$ws.Range("J98").Value = 15000
$ws.Range("L98").Value = 15000
$ws.Range("N98").Value = -20990
# row 102
$ws.Range("H102").Value = 4207.273
$ws.Range("I102").Value = 897.1429000000001
$ws.Range("K102").Value = 897.1429000000001
$ws.Range("M102").Value = 724.8570999999999
# row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = ""
$ws.Range("N114").Value = 0
# row 116
$ws.Range("H116").Value = 3550.5
$ws.Range("I116").Value = 2889.5
$ws.Range("K116").Value = 2889.5
$ws.Range("M116").Value = -595.5
# row 122
$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50
# row 136
$ws.Range("H136").Value = 2864
$ws.Range("I136").Value = 2436.8
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7310.400000000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4760.400000000001
$ws.Range("N136").Value = -20100

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 3550.5
$ws.Range("I3").Value = 2889.5
$ws.Range("K3").Value = 2889.5
$ws.Range("M3").Value = -2775.5
# row 22
$ws.Range("H22").Value = 216.66667
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27
# row 86
$ws.Range("H86").Value = 3315.8635
$ws.Range("J86").Value = 5743.7144
$ws.Range("L86").Value = 5743.7144
$ws.Range("N86").Value = -7989.7144
# row 89
$ws.Range("H89").Value = 3315.8635
$ws.Range("J89").Value = 5743.7144
$ws.Range("L89").Value = 28718.572
$ws.Range("N89").Value = -39950.572
# row 94
$ws.Range("H94").Value = 240.35294
$ws.Range("I94").Value = 216.46153
$ws.Range("J94").Value = 318
$ws.Range("K94").Value = 216.46153
$ws.Range("L94").Value = 318
$ws.Range("M94").Value = 234.53847
$ws.Range("N94").Value = -1220
# row 99
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
# row 107
$ws.Range("H107").Value = 3825.7273
$ws.Range("I107").Value = 894.6923
$ws.Range("J107").Value = 8059.4443
$ws.Range("K107").Value = 894.6923
$ws.Range("L107").Value = 8059.4443
$ws.Range("M107").Value = 1025.3077
$ws.Range("N107").Value = -11899.4443

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 6550.316
$ws.Range("I31").Value = 4061.182
$ws.Range("J31").Value = 7564.407
$ws.Range("K31").Value = 4061.182
$ws.Range("L31").Value = 7564.407
$ws.Range("M31").Value = -3766.182
$ws.Range("N31").Value = -8154.407
# row 34
$ws.Range("H34").Value = 6550.316
$ws.Range("I34").Value = 4061.182
$ws.Range("J34").Value = 7564.407
$ws.Range("K34").Value = 4061.182
$ws.Range("L34").Value = 7564.407
$ws.Range("M34").Value = -3859.182
$ws.Range("N34").Value = -7968.407
# row 47
$ws.Range("H47").Value = 64
$ws.Range("I47").Value = 64
$ws.Range("K47").Value = 64
$ws.Range("M47").Value = 502
# row 48
$ws.Range("H48").Value = 49999
$ws.Range("J48").Value = 49999
$ws.Range("L48").Value = 49999
$ws.Range("N48").Value = -50951
# row 58
$ws.Range("H58").Value = 4331.5557
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1797
# row 62
$ws.Range("H62").Value = 4937.5
$ws.Range("J62").Value = 4937.5
$ws.Range("L62").Value = 4937.5
$ws.Range("N62").Value = -6185.5
# row 65
$ws.Range("H65").Value = 4937.5
$ws.Range("J65").Value = 4937.5
$ws.Range("L65").Value = 24687.5
$ws.Range("N65").Value = -30927.5
# row 74
$ws.Range("H74").Value = 58439.11
$ws.Range("J74").Value = 58439.11
$ws.Range("L74").Value = 58439.11
$ws.Range("N74").Value = -60187.11
# row 77
$ws.Range("H77").Value = 58439.11
$ws.Range("J77").Value = 58439.11
$ws.Range("L77").Value = 175317.33
$ws.Range("N77").Value = -184053.33
# row 94
$ws.Range("H94").Value = 8103.25
$ws.Range("I94").Value = 1799
$ws.Range("J94").Value = 10204.667
$ws.Range("K94").Value = 1799
$ws.Range("L94").Value = 10204.667
$ws.Range("M94").Value = -1348
$ws.Range("N94").Value = -11106.667
# row 136
$ws.Range("H136").Value = 4331.5557
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450
# row 141
$ws.Range("H141").Value = 86239.8
$ws.Range("J141").Value = 86239.8
$ws.Range("L141").Value = 86239.8
$ws.Range("N141").Value = -96599.8

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 67.888885
$ws.Range("I2").Value = 107.75
$ws.Range("J2").Value = 36
$ws.Range("K2").Value = 646.5
$ws.Range("L2").Value = 216
$ws.Range("M2").Value = -533.5
$ws.Range("N2").Value = -442
# row 11
$ws.Range("H11").Value = 126.38461
$ws.Range("I11").Value = 48.833332
$ws.Range("K11").Value = 146.499996
$ws.Range("M11").Value = -6.49999600000001
# row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""
# row 33
$ws.Range("H33").Value = 67.8
$ws.Range("I33").Value = 64
$ws.Range("J33").Value = 71.125
$ws.Range("K33").Value = 384
$ws.Range("L33").Value = 426.75
$ws.Range("M33").Value = -101
$ws.Range("N33").Value = -992.75
# row 81
$ws.Range("H81").Value = 3400
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 3400
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = ""
$ws.Range("M81").Value = 10200
$ws.Range("N81").Value = -12446
# row 84
$ws.Range("H84").Value = 3400
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3400
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = ""
$ws.Range("M84").Value = 30600
$ws.Range("N84").Value = -41832
# row 110
$ws.Range("H110").Value = 3501
$ws.Range("I110").Value = 3501.25
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 10503.75
$ws.Range("L110").Value = 10500
$ws.Range("M110").Value = -6413.75
$ws.Range("N110").Value = -18680
# row 121
$ws.Range("H121").Value = 369
$ws.Range("J121").Value = 999.5
$ws.Range("L121").Value = 2998.5
$ws.Range("N121").Value = -5618.5
# row 128
$ws.Range("H128").Value = 533328.2
$ws.Range("I128").Value = 533328.2
$ws.Range("K128").Value = 1599984.6
$ws.Range("M128").Value = -1595004.6

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = ""
$ws.Range("N26").Value = 0
# row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = ""
$ws.Range("N50").Value = 0
# row 52
$ws.Range("H52").Value = 55000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 55000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = ""
$ws.Range("M52").Value = 55000
$ws.Range("N52").Value = -55518
# row 80
$ws.Range("H80").Value = 5492.5
$ws.Range("I80").Value = 5492.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5492.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -4494.5
# row 83
$ws.Range("H83").Value = 5492.5
$ws.Range("I83").Value = 5492.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 27462.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -22470.5
# row 126
$ws.Range("H126").Value = 6100.5386
$ws.Range("I126").Value = 5644.778
$ws.Range("J126").Value = 7126
$ws.Range("K126").Value = 16934.334
$ws.Range("L126").Value = 21378
$ws.Range("M126").Value = -14464.334
$ws.Range("N126").Value = -26318
# row 132
$ws.Range("H132").Value = 5430.35
$ws.Range("I132").Value = 5451.067
$ws.Range("J132").Value = 5368.2
$ws.Range("K132").Value = 16353.201
$ws.Range("L132").Value = 16104.6
$ws.Range("M132").Value = -13823.201
$ws.Range("N132").Value = -21164.6
# row 134
$ws.Range("H134").Value = 97591.164
$ws.Range("J134").Value = 97591.164
$ws.Range("L134").Value = 292773.492
$ws.Range("N134").Value = -297843.492

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# row 5
$ws.Range("H5").Value = 61110
$ws.Range("J5").Value = 61110
$ws.Range("L5").Value = 61110
$ws.Range("N5").Value = -61336
# row 7
$ws.Range("H7").Value = 5211.0835
$ws.Range("I7").Value = 3734
$ws.Range("J7").Value = 7279
$ws.Range("K7").Value = 3734
$ws.Range("L7").Value = 7279
$ws.Range("M7").Value = -3622
$ws.Range("N7").Value = -7503
# row 16
$ws.Range("H16").Value = 193.2
$ws.Range("I16").Value = 193.2
$ws.Range("K16").Value = 193.2
$ws.Range("M16").Value = -23.19999999999999
# row 20
$ws.Range("H20").Value = 507499.5
$ws.Range("J20").Value = 507499.5
$ws.Range("L20").Value = 507499.5
$ws.Range("N20").Value = -507951.5
# row 22
$ws.Range("H22").Value = 919.1
$ws.Range("J22").Value = 849.2
$ws.Range("L22").Value = 849.2
$ws.Range("N22").Value = -1439.2
# row 27
$ws.Range("H27").Value = 919.1
$ws.Range("J27").Value = 849.2
$ws.Range("L27").Value = 849.2
$ws.Range("N27").Value = -1063.2
# row 46
$ws.Range("H46").Value = 5061.1
$ws.Range("J46").Value = 5199.6665
$ws.Range("L46").Value = 5199.6665
$ws.Range("N46").Value = -5575.6665
# row 50
$ws.Range("H50").Value = 73074
$ws.Range("I50").Value = 73074
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 73074
$ws.Range("L50").Value = ""
$ws.Range("N50").Value = 0
$ws.Range("M50").Value = -72437
# row 61
$ws.Range("H61").Value = 4099.7856
$ws.Range("I61").Value = 1924.875
$ws.Range("J61").Value = 6999.6665
$ws.Range("K61").Value = 1924.875
$ws.Range("L61").Value = 6999.6665
$ws.Range("M61").Value = -1722.875
$ws.Range("N61").Value = -7403.6665
# row 68
$ws.Range("H68").Value = 7968.625
$ws.Range("I68").Value = 4916.3335
$ws.Range("J68").Value = 9800
$ws.Range("K68").Value = 4916.3335
$ws.Range("L68").Value = 9800
$ws.Range("M68").Value = -4167.3335
$ws.Range("N68").Value = -11298
# row 71
$ws.Range("H71").Value = 7968.625
$ws.Range("I71").Value = 4916.3335
$ws.Range("J71").Value = 9800
$ws.Range("K71").Value = 24581.6675
$ws.Range("L71").Value = 49000
$ws.Range("M71").Value = -20837.6675
$ws.Range("N71").Value = -56488
# row 93
$ws.Range("H93").Value = 1588.6666
$ws.Range("I93").Value = 1485.6
$ws.Range("K93").Value = 1485.6
$ws.Range("M93").Value = -237.5999999999999
# row 113
$ws.Range("H113").Value = 4099.7856
$ws.Range("I113").Value = 1924.875
$ws.Range("J113").Value = 6999.6665
$ws.Range("K113").Value = 1924.875
$ws.Range("L113").Value = 6999.6665
$ws.Range("M113").Value = 245.125
$ws.Range("N113").Value = -11339.6665
# row 122
$ws.Range("H122").Value = 3001.1667
$ws.Range("I122").Value = 2934.8
$ws.Range("K122").Value = 8804.400000000001
$ws.Range("M122").Value = -6354.400000000001
# row 126
$ws.Range("H126").Value = 5211.0835
$ws.Range("I126").Value = 3734
$ws.Range("J126").Value = 7279
$ws.Range("K126").Value = 11202
$ws.Range("L126").Value = 21837
$ws.Range("M126").Value = -8732
$ws.Range("N126").Value = -26777
# row 132
$ws.Range("H132").Value = 3759
$ws.Range("I132").Value = 3498.5
$ws.Range("J132").Value = 3932.6667
$ws.Range("K132").Value = 10495.5
$ws.Range("L132").Value = 11798.0001
$ws.Range("M132").Value = -7965.5
$ws.Range("N132").Value = -16858.0001

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = ""
$ws.Range("N39").Value = 0
# row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = ""
$ws.Range("N42").Value = 0
# row 43
$ws.Range("H43").Value = 22333.334
$ws.Range("J43").Value = 22333.334
$ws.Range("L43").Value = 22333.334
$ws.Range("N43").Value = -22631.334
# row 51
$ws.Range("H51").Value = 13750
$ws.Range("I51").Value = 13750
$ws.Range("K51").Value = 13750
$ws.Range("M51").Value = -13240
# row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = ""
$ws.Range("N80").Value = 0
# row 81
$ws.Range("H81").Value = 1333.3334
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 4000
$ws.Range("N81").Value = -6122
# row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = ""
$ws.Range("N82").Value = 0
# row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = ""
$ws.Range("N83").Value = 0
# row 84
$ws.Range("H84").Value = 1333.3334
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 20000
$ws.Range("N84").Value = -30608
# row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = ""
$ws.Range("N85").Value = 0
# row 88
$ws.Range("H88").Value = 27500
$ws.Range("J88").Value = 27500
$ws.Range("L88").Value = 27500
$ws.Range("N88").Value = -28312
# row 91
$ws.Range("H91").Value = 27500
$ws.Range("J91").Value = 27500
$ws.Range("L91").Value = 27500
$ws.Range("N91").Value = -30308
# row 96
$ws.Range("H96").Value = 1576.8182
$ws.Range("I96").Value = 1571.2858
$ws.Range("J96").Value = 1586.5
$ws.Range("K96").Value = 1571.2858
$ws.Range("L96").Value = 1586.5
$ws.Range("M96").Value = -198.2858000000001
$ws.Range("N96").Value = -4332.5
# row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = ""
$ws.Range("N135").Value = 0
